$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update "Pagos" (F) and "Inscrições homologadas" (H) counts for the affected rows.
$ws.Range("F8").Value = 13
$ws.Range("H8").Value = 14

$ws.Range("F14").Value = 2
$ws.Range("H14").Value = 3

$ws.Range("F16").Value = 6
$ws.Range("H16").Value = 10

$ws.Range("F27").Value = 17
$ws.Range("H27").Value = 21

$ws.Range("F35").Value = 11
$ws.Range("H35").Value = 12

$ws.Range("F37").Value = 47
$ws.Range("H37").Value = 59

$ws.Range("F38").Value = 33
$ws.Range("H38").Value = 53

$ws.Range("F61").Value = 23
$ws.Range("H61").Value = 33

$ws.Range("F68").Value = 12
$ws.Range("H68").Value = 16
